# Apply the authored edit described by the commit diff:
#  - shrink the saved window size in bookViews
#  - replace the "Solana" tx sample row with a single "XRP" tx row
#  - clear out the other now-unused sample tx rows (3-7), but keep their
#    existing cell style
#  - move the active selection from C3 to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Saved window size (bookViews/workbookView) ---------------------------
$win = $wb.Windows.Item(1)
$win.Width  = 30240
$win.Height = 11500

# --- Row 2: swap the Solana sample for an XRP one --------------------------
$ws.Range("B2").Value = "XRP"
$ws.Range("C2").Value = "2EB4913CE256D6ABFAC9CED396E67F9DAD3DF8D92FB9B9B91F2D6481F6EF34E3"

# The new values were typed fresh (no distributed-alignment style carried
# over), so strip the style off B2:C2 by pasting in formats from a cell
# that never had one, then clean that helper cell back up.
$ws.Range("D1").Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)
$ws.Range("D1").Delete()

# Row 2 no longer needs its custom row height.
$ws.Rows.Item(2).AutoFit()

# --- Rows 3-7: drop the remaining sample tx rows ---------------------------
$ws.Range("A3:C7").ClearContents()
for ($r = 3; $r -le 7; $r++) {
  $ws.Rows.Item($r).AutoFit()
}

# --- Selection moves from C3 to C2 -----------------------------------------
$ws.Range("C2").Select() | Out-Null
